# Actualización automática del mapa: agrega los nuevos casos (filas 66-69)
# a la hoja "NEW" del mapa interactivo.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

# Caso, F. De Reclamo, Direccion, Comuna, OT, Proveedor Asignado, Estado,
# Observaciones, Attachments, Tipo de tarea, Equipo, Tipo de Elemento
$newCases = @(
    @("6943", "8/14/2025", "SUPERI 1459", "13", "808972965", "NEW", "Pendiente", "Poste con base quebrada ver si es posible desmonte", 1, "Desmonte", "Sin equipos", "Poste"),
    @("6944", "8/14/2025", "RAVIGNANI, EMILIO, DR. 2040", "14", "808972970", "NEW", "Pendiente", "Picada", 1, "Cambio", "Sin equipos", "Terminal"),
    @("6969", "8/14/2025", "CIUDAD DE LA PAZ 295", "14", "808972995", "NEW", "Pendiente", "Cambiar", 1, "Cambio", "Sin equipos", "Pasante"),
    @("6971", "8/14/2025", "MAURE 1594", "14", "808973001", "NEW", "Pendiente", "Ver de traspasar a telecentro y  desmontar ver con inspector", 1, "Cambio", "Sin equipos", "Pasante")
)

$startRow = 66
# Columns that must stay text even though they look numeric/date-like
# (Caso, F. De Reclamo, Comuna, OT) — matches the rest of the sheet where
# these columns are stored as plain text, not numbers/dates.
$textColumns = @(1, 2, 4, 5)

for ($i = 0; $i -lt $newCases.Count; $i++) {
    $r = $startRow + $i
    $values = $newCases[$i]

    for ($c = 1; $c -le $values.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($textColumns -contains $c) {
            # Force text entry so "6943", "8/14/2025", "13", "808972965"
            # are not silently coerced into numbers/dates, then drop the
            # temporary text format so the cell keeps the default style.
            $cell.NumberFormat = "@"
            $cell.Value = $values[$c - 1]
            $cell.ClearFormats()
        } else {
            $cell.Value = $values[$c - 1]
        }
    }
    # Coordenada_X, Coordenada_Y, Operacion, Zona are left blank for these
    # newly reported cases (no geocoding/zone assigned yet).
}
